$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the raw (unrounded) input values that feed the ROUND()/SUM() formulas.
# Dependent cells (D/E/F columns) recalculate automatically.
$ws.Range("H5").Value = 90.219406127929688
$ws.Range("I6").Value = 9.7805938720703125
$ws.Range("H8").Value = 4.3786153793334961
$ws.Range("I8").Value = 5.2050833702087402
$ws.Range("H9").Value = 6.8983349800109863
$ws.Range("I9").Value = 7.7633652687072754
$ws.Range("H10").Value = 93.805458068847656
$ws.Range("I10").Value = 37.183757781982422
$ws.Range("J10").Value = 87.40167236328125
